# Update "想去人数" (interested-people count) figures that changed between
# crawl runs for two of the duplicated event rows, on both the "展览"
# sheet and the consolidated "全部类型" sheet.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F5").Value = 2993
$wsExhibit.Range("F6").Value = 302
$wsExhibit.Range("F7").Value = 403

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 2993
$wsAll.Range("F6").Value = 302
$wsAll.Range("F9").Value = 403
